$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new single supplier example
$ws.Range("A2").Value = "TO"
$ws.Range("B2").Value = "The Originote"
$ws.Range("C2").Value = "Jl.Bali Panjang 31"

# Remove the now unused third row (PGN / Piggeon / Jl.Bromo No 11)
$ws.Rows.Item(3).Delete()

# Update the selected cell/view as in the saved file
$ws.Range("G10").Select()
